# Fixed Employee Count Bug
# Updates "Average Comp" / "Total Employees" figures on the per-year detail
# sheet ("Sheet") and mirrors the corresponding values on the "WPI" summary
# sheet, then removes the duplicated header row that used to sit between
# the "Total Reported Employees" label row and the real "Year/Name/..."
# header row (its counts now live directly on the label row instead).

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item(1)   # "Sheet" - per-year raw data
$wsWPI  = $wb.Worksheets.Item(2)   # "WPI"   - summary/report sheet

# New Average Comp (text, formatted like the existing cells) and
# Total Employees (numeric) values, one row per fiscal year 2022 -> 2014.
$rows = @(
    @{ Row = 2;  AvgComp = '$449,394.53'; TotalEmp = 19 },
    @{ Row = 3;  AvgComp = '$397,600.67'; TotalEmp = 18 },
    @{ Row = 4;  AvgComp = '$496,126.71'; TotalEmp = 17 },
    @{ Row = 5;  AvgComp = '$437,979.28'; TotalEmp = 18 },
    @{ Row = 6;  AvgComp = '$407,612.42'; TotalEmp = 19 },
    @{ Row = 7;  AvgComp = '$439,114.13'; TotalEmp = 15 },
    @{ Row = 8;  AvgComp = '$385,720.19'; TotalEmp = 16 },
    @{ Row = 9;  AvgComp = '$392,763.25'; TotalEmp = 16 },
    @{ Row = 10; AvgComp = '$386,920.81'; TotalEmp = 16 }
)

foreach ($r in $rows) {
    $cell = $wsData.Cells.Item($r.Row, 6)   # column F = Average Comp
    $cell.NumberFormat = "@"
    $cell.Value = $r.AvgComp
    $wsData.Cells.Item($r.Row, 8).Value = $r.TotalEmp   # column H = Total Employees
}

# Mirror the "Average Comp Per Reported Employee" row (row 7, columns B:J)
# on the WPI summary sheet.
$wpiCols = @("B", "C", "D", "E", "F", "G", "H", "I", "J")
for ($i = 0; $i -lt $rows.Count; $i++) {
    $col = $wpiCols[$i]
    $cell = $wsWPI.Range($col + "7")
    $cell.NumberFormat = "@"
    $cell.Value = $rows[$i].AvgComp
}

# Populate the "Total Reported Employees" row (row 9) with the per-year
# counts, directly under row 9's label in column A.
for ($i = 0; $i -lt $rows.Count; $i++) {
    $col = $wpiCols[$i]
    $wsWPI.Range($col + "9").Value = $rows[$i].TotalEmp
}

# Remove the now-redundant duplicate header row (old row 10: Name/Title/
# Title Group/Base Compensation/Other Comp/Total Comp) that duplicated the
# real header in row 11. Clearing it collapses/removes the row entirely.
$wsWPI.Range("B10:G10").ClearContents()
